$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-16 all hold the same date serial value (46081)
# that needs to be bumped forward by one day to 46082.
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
